$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -10
$ws.Range("F5").Value = -7
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = -2
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -7
$ws.Range("F22").Value = -11
$ws.Range("F24").Value = 8
